$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column B to text format so zero-padded IDs keep their leading zeros
$ws.Range("B2:B71").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "│¬└╠"
$ws.Cells.Item(2, 2).Value = "0001"
$ws.Cells.Item(3, 1).Value = "│¬┤⌐┤┘"
$ws.Cells.Item(3, 2).Value = "0002"
$ws.Cells.Item(4, 1).Value = "│¬┼╕│¬┤┘"
$ws.Cells.Item(4, 2).Value = "0003"
$ws.Cells.Item(5, 1).Value = "│¬░í┤┘"
$ws.Cells.Item(5, 2).Value = "0004"
$ws.Cells.Item(6, 1).Value = "│δ╖┬"
$ws.Cells.Item(6, 2).Value = "0005"
$ws.Cells.Item(7, 1).Value = "│ε┤┘"
$ws.Cells.Item(7, 2).Value = "0006"
$ws.Cells.Item(8, 1).Value = "│ε╢≤┤┘"
$ws.Cells.Item(8, 2).Value = "0007"
$ws.Cells.Item(9, 1).Value = "│⌡─í┤┘"
$ws.Cells.Item(9, 2).Value = "0008"
$ws.Cells.Item(10, 1).Value = "│╖┤┘"
$ws.Cells.Item(10, 2).Value = "0009"
$ws.Cells.Item(11, 1).Value = "│╗└╧"
$ws.Cells.Item(11, 2).Value = "0010"
$ws.Cells.Item(12, 1).Value = "│╗╕«┤┘"
$ws.Cells.Item(12, 2).Value = "0011"
$ws.Cells.Item(13, 1).Value = "│╗╖┴│⌡┤┘"
$ws.Cells.Item(13, 2).Value = "0012"
$ws.Cells.Item(14, 1).Value = "│╤┤┘"
$ws.Cells.Item(14, 2).Value = "0013"
$ws.Cells.Item(15, 1).Value = "│╤╛ε┐└┤┘"
$ws.Cells.Item(15, 2).Value = "0014"
$ws.Cells.Item(16, 1).Value = "│╩"
$ws.Cells.Item(16, 2).Value = "0015"
$ws.Cells.Item(17, 1).Value = "│▓└┌"
$ws.Cells.Item(17, 2).Value = "0016"
$ws.Cells.Item(18, 1).Value = "│▓┤┘"
$ws.Cells.Item(18, 2).Value = "0017"
$ws.Cells.Item(19, 1).Value = "│▓▒Γ┤┘"
$ws.Cells.Item(19, 2).Value = "0018"
$ws.Cells.Item(20, 1).Value = "┐°╟╧┤┘"
$ws.Cells.Item(20, 2).Value = "0019"
$ws.Cells.Item(21, 1).Value = "┐⌠┤┘"
$ws.Cells.Item(21, 2).Value = "0020"
$ws.Cells.Item(22, 1).Value = "┐└┤├"
$ws.Cells.Item(22, 2).Value = "0021"
$ws.Cells.Item(23, 1).Value = "┐┬╡╡"
$ws.Cells.Item(23, 2).Value = "0022"
$ws.Cells.Item(24, 1).Value = "┐╓"
$ws.Cells.Item(24, 2).Value = "0023"
$ws.Cells.Item(25, 1).Value = "└╘┐°"
$ws.Cells.Item(25, 2).Value = "0024"
$ws.Cells.Item(26, 1).Value = "└╠╕º"
$ws.Cells.Item(26, 2).Value = "0025"
$ws.Cells.Item(27, 1).Value = "└╠╟╪"
$ws.Cells.Item(27, 2).Value = "0026"
$ws.Cells.Item(28, 1).Value = "└╥╛ε╣÷╕«┤┘"
$ws.Cells.Item(28, 2).Value = "0027"
$ws.Cells.Item(29, 1).Value = "└╧║╬╖»"
$ws.Cells.Item(29, 2).Value = "0028"
$ws.Cells.Item(30, 1).Value = "└╪┤┘"
$ws.Cells.Item(30, 2).Value = "0029"
$ws.Cells.Item(31, 1).Value = "┤÷║╨┐í"
$ws.Cells.Item(31, 2).Value = "0030"
$ws.Cells.Item(32, 1).Value = "┤δ╚¡"
$ws.Cells.Item(32, 2).Value = "0031"
$ws.Cells.Item(33, 1).Value = "┤⌐▒╕"
$ws.Cells.Item(33, 2).Value = "0032"
$ws.Cells.Item(34, 1).Value = "┤┘└╜"
$ws.Cells.Item(34, 2).Value = "0033"
$ws.Cells.Item(35, 1).Value = "┤┘╕«"
$ws.Cells.Item(35, 2).Value = "0034"
$ws.Cells.Item(36, 1).Value = "┤├╛ε│¬┤┘"
$ws.Cells.Item(36, 2).Value = "0035"
$ws.Cells.Item(37, 1).Value = "┤▌┤┘"
$ws.Cells.Item(37, 2).Value = "0036"
$ws.Cells.Item(38, 1).Value = "╛α╝╙"
$ws.Cells.Item(38, 2).Value = "0037"
$ws.Cells.Item(39, 1).Value = "╛ε┴÷╖┤┤┘"
$ws.Cells.Item(39, 2).Value = "0038"
$ws.Cells.Item(40, 1).Value = "╛ε╕░└╠"
$ws.Cells.Item(40, 2).Value = "0039"
$ws.Cells.Item(41, 1).Value = "╛ε╖╞┤┘"
$ws.Cells.Item(41, 2).Value = "0040"
$ws.Cells.Item(42, 1).Value = "╛ε╢╗░╘"
$ws.Cells.Item(42, 2).Value = "0041"
$ws.Cells.Item(43, 1).Value = "╛╚│τ╟╧╝╝┐Σ"
$ws.Cells.Item(43, 2).Value = "0042"
$ws.Cells.Item(44, 1).Value = "╛╚┼╕▒⌡┤┘"
$ws.Cells.Item(44, 2).Value = "0043"
$ws.Cells.Item(45, 1).Value = "╛╚╜╔"
$ws.Cells.Item(45, 2).Value = "0044"
$ws.Cells.Item(46, 1).Value = "╛╞│ó┤┘"
$ws.Cells.Item(46, 2).Value = "0045"
$ws.Cells.Item(47, 1).Value = "╛╞╕º┤Σ┤┘"
$ws.Cells.Item(47, 2).Value = "0046"
$ws.Cells.Item(48, 1).Value = "╛╦┤┘"
$ws.Cells.Item(48, 2).Value = "0047"
$ws.Cells.Item(49, 1).Value = "╛╦╖┴┴╓┤┘"
$ws.Cells.Item(49, 2).Value = "0048"
$ws.Cells.Item(50, 1).Value = "╛╦╛╞┬≈╕«┤┘"
$ws.Cells.Item(50, 2).Value = "0049"
$ws.Cells.Item(51, 1).Value = "╛╦╛╞╝¡╟╧┤┘"
$ws.Cells.Item(51, 2).Value = "0050"
$ws.Cells.Item(52, 1).Value = "╡╡┐≥"
$ws.Cells.Item(52, 2).Value = "0051"
$ws.Cells.Item(53, 1).Value = "╡╢╝¡"
$ws.Cells.Item(53, 2).Value = "0052"
$ws.Cells.Item(54, 1).Value = "░°║╬"
$ws.Cells.Item(54, 2).Value = "0053"
$ws.Cells.Item(55, 1).Value = "░µ╟Φ"
$ws.Cells.Item(55, 2).Value = "0054"
$ws.Cells.Item(56, 1).Value = "░¿╗τ"
$ws.Cells.Item(56, 2).Value = "0055"
$ws.Cells.Item(57, 1).Value = "░í─í"
$ws.Cells.Item(57, 2).Value = "0056"
$ws.Cells.Item(58, 1).Value = "░ß┴ñ"
$ws.Cells.Item(58, 2).Value = "0057"
$ws.Cells.Item(59, 1).Value = "░φ┼δ"
$ws.Cells.Item(59, 2).Value = "0058"
$ws.Cells.Item(60, 1).Value = "░Φ╚╣"
$ws.Cells.Item(60, 2).Value = "0059"
$ws.Cells.Item(61, 1).Value = "░φ╣╬"
$ws.Cells.Item(61, 2).Value = "0060"
$ws.Cells.Item(62, 1).Value = "░ⁿ╜╔"
$ws.Cells.Item(62, 2).Value = "0061"
$ws.Cells.Item(63, 1).Value = "░╞┴ñ"
$ws.Cells.Item(63, 2).Value = "0062"
$ws.Cells.Item(64, 1).Value = "░╟░¡"
$ws.Cells.Item(64, 2).Value = "0063"
$ws.Cells.Item(65, 1).Value = "▒Γ└√"
$ws.Cells.Item(65, 2).Value = "0064"
$ws.Cells.Item(66, 1).Value = "▒Γ┤δ"
$ws.Cells.Item(66, 2).Value = "0065"
$ws.Cells.Item(67, 1).Value = "▒Γ║╨"
$ws.Cells.Item(67, 2).Value = "0066"
$ws.Cells.Item(68, 1).Value = "▒Γ╗▌"
$ws.Cells.Item(68, 2).Value = "0067"
$ws.Cells.Item(69, 1).Value = "▒Γ╚╕"
$ws.Cells.Item(69, 2).Value = "0068"
$ws.Cells.Item(70, 1).Value = "▒Γ╛∩"
$ws.Cells.Item(70, 2).Value = "0069"
$ws.Cells.Item(71, 1).Value = "▒│└░"
$ws.Cells.Item(71, 2).Value = "0070"
